# Sprint_Backlog.xlsx update: reorder/insert rows on the "Backlog" sheet
# - Row 23 gets a brand-new task "Close chat if its deleted"
# - Row 26 becomes "Show newest chat info" (renamed from "Show newest info")
# - Rows 20-44 are rewritten to reflect the new backlog ordering
# - sheetView scroll position / selection updated to match

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Backlog")

# New contents for column A (priority) / column B (task name), rows 20..44
$rows = @(
  @(20, 2, "Show message info"),
  @(21, 1, "Remove your message"),
  @(22, 2, "Show if someone is writing"),
  @(23, 2, "Close chat if its deleted"),
  @(24, 2, "Show online persons"),
  @(25, 2, "Show chat status (online P.)"),
  @(26, 2, "Show newest chat info"),
  @(27, 2, "Play song"),
  @(28, 2, "Save song"),
  @(29, 2, "Get songs info"),
  @(30, 2, "Update song"),
  @(31, 2, "Remove song"),
  @(32, 3, "Change time displaying"),
  @(33, 2, "Notifications"),
  @(34, 3, "Multiplayer game"),
  @(35, 3, "Theme"),
  @(36, 3, "Like"),
  @(37, 3, "Add Friend"),
  @(38, 3, "Remove Friend"),
  @(39, 3, "Manage group"),
  @(40, 1, "Dedicated client"),
  @(41, 2, "Web client"),
  @(42, 3, "News area"),
  @(43, 3, "Sound equalizer"),
  @(44, 2, "Join as group")
)

foreach ($row in $rows) {
    $r = $row[0]
    $priority = $row[1]
    $name = $row[2]
    $ws.Cells.Item($r, 1).Value = $priority
    $ws.Cells.Item($r, 2).Value = $name
}

# Update the view: scroll position and active selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B26").Select()
